$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="60÷2=30, 0";  New="51÷7=7, 2"},
    @{Row=1;  Col=2; Old="91÷7=13, 0"; New="11÷3=3, 2"},
    @{Row=1;  Col=3; Old="99÷8=12, 3"; New="87÷9=9, 6"},
    @{Row=1;  Col=4; Old="12÷8=1, 4";  New="63÷4=15, 3"},
    @{Row=1;  Col=5; Old="95÷2=47, 1"; New="91÷3=30, 1"},
    @{Row=5;  Col=1; Old="17÷5=3, 2";  New="66÷8=8, 2"},
    @{Row=5;  Col=2; Old="56÷4=14, 0"; New="79÷5=15, 4"},
    @{Row=5;  Col=3; Old="84÷2=42, 0"; New="15÷8=1, 7"},
    @{Row=5;  Col=4; Old="58÷8=7, 2";  New="10÷6=1, 4"},
    @{Row=5;  Col=5; Old="79÷2=39, 1"; New="94÷3=31, 1"},
    @{Row=9;  Col=1; Old="65÷9=7, 2";  New="42÷8=5, 2"},
    @{Row=9;  Col=2; Old="41÷4=10, 1"; New="52÷9=5, 7"},
    @{Row=9;  Col=3; Old="17÷5=3, 2";  New="35÷7=5, 0"},
    @{Row=9;  Col=4; Old="40÷5=8, 0";  New="82÷3=27, 1"},
    @{Row=9;  Col=5; Old="51÷2=25, 1"; New="32÷8=4, 0"},
    @{Row=13; Col=1; Old="47÷7=6, 5";  New="58÷5=11, 3"},
    @{Row=13; Col=2; Old="85÷4=21, 1"; New="75÷2=37, 1"},
    @{Row=13; Col=3; Old="29÷2=14, 1"; New="82÷2=41, 0"},
    @{Row=13; Col=4; Old="80÷5=16, 0"; New="20÷9=2, 2"},
    @{Row=13; Col=5; Old="97÷7=13, 6"; New="52÷4=13, 0"},
    @{Row=17; Col=1; Old="14÷2=7, 0";  New="71÷7=10, 1"},
    @{Row=17; Col=2; Old="78÷7=11, 1"; New="77÷2=38, 1"},
    @{Row=17; Col=3; Old="35÷8=4, 3";  New="35÷4=8, 3"},
    @{Row=17; Col=4; Old="85÷7=12, 1"; New="39÷3=13, 0"},
    @{Row=17; Col=5; Old="15÷2=7, 1";  New="94÷4=23, 2"}
)

foreach ($r in $replacements) {
    $cell = $tbl.Cell($r.Row, $r.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
